# Weekly fruit/vegetable price update: insert a new week's worth of rows
# (2 rows: "Primera" / "Segunda" quality) at the top of the data table,
# pushing all existing rows down by 2 (dimension grows from R370 to R372).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 345-346; everything currently at/after row 345
# (through 370) shifts down to 347-372.
$ws.Rows("345:346").Insert()

# New row 345: Betarraga, "Primera" quality, week of 44578
$ws.Cells.Item(345, 1).Value  = 9
$ws.Cells.Item(345, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(345, 3).Value  = "Metropolitana"
$ws.Cells.Item(345, 4).Value  = 44578
$ws.Cells.Item(345, 5).Value  = 13
$ws.Cells.Item(345, 6).Value  = 100114014
$ws.Cells.Item(345, 7).Value  = "Betarraga"
$ws.Cells.Item(345, 8).Value  = "Sin especificar"
$ws.Cells.Item(345, 9).Value  = "Primera"
$ws.Cells.Item(345, 10).Value = 4300
$ws.Cells.Item(345, 11).Value = 100
$ws.Cells.Item(345, 12).Value = 120
$ws.Cells.Item(345, 13).Value = 110
$ws.Cells.Item(345, 14).Value = "$/unidad"
$ws.Cells.Item(345, 15).Value = "Región Metropolitana"
$ws.Cells.Item(345, 16).Value = 110
$ws.Cells.Item(345, 17).Value = 1
$ws.Cells.Item(345, 18).Value = "Hortaliza"

# New row 346: Betarraga, "Segunda" quality, week of 44578
$ws.Cells.Item(346, 1).Value  = 9
$ws.Cells.Item(346, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(346, 3).Value  = "Metropolitana"
$ws.Cells.Item(346, 4).Value  = 44578
$ws.Cells.Item(346, 5).Value  = 13
$ws.Cells.Item(346, 6).Value  = 100114014
$ws.Cells.Item(346, 7).Value  = "Betarraga"
$ws.Cells.Item(346, 8).Value  = "Sin especificar"
$ws.Cells.Item(346, 9).Value  = "Segunda"
$ws.Cells.Item(346, 10).Value = 1960
$ws.Cells.Item(346, 11).Value = 70
$ws.Cells.Item(346, 12).Value = 80
$ws.Cells.Item(346, 13).Value = 75
$ws.Cells.Item(346, 14).Value = "$/unidad"
$ws.Cells.Item(346, 15).Value = "Región Metropolitana"
$ws.Cells.Item(346, 16).Value = 75
$ws.Cells.Item(346, 17).Value = 1
$ws.Cells.Item(346, 18).Value = "Hortaliza"
